# Reorder the two worksheets so "review_info" precedes "hotel_info",
# and insert a new "State" column (value "Louisiana") right after
# "Hotel_Name" in the hotel_info sheet.

$wb = $excel.ActiveWorkbook

$hotelWs  = $wb.Worksheets.Item("hotel_info")
$reviewWs = $wb.Worksheets.Item("review_info")

# Insert a new column before the current City column (column C) in hotel_info,
# shifting City/Zip/TA_ReviewURL/Tripadvisor_Hotel_Name/English_Reviews_num/
# Local_Rank/Total_Reviews_num one column to the right.
$hotelWs.Range("C:C").Insert()

$hotelWs.Range("C1").Value = "State"
$hotelWs.Range("C2").Value = "Louisiana"

# Make review_info the first tab, hotel_info the second.
$reviewWs.Move($hotelWs)
